{"js": "// Update the date heading at the top of the document.\nconst firstParagraph = context.document.body.paragraphs.getFirst();\nconst dateRange = firstParagraph.getRange();\ndateRange.insertText(\"2026-02-14 Saturday\", Word.InsertLocation.replace);\n\n// Update each three-digit x one-digit multiplication answer in the table.\n// The table has data only in rows 0, 4, 9, 14, 19 (0-indexed; 5 columns\n// each); replacements are applied in row-major order so duplicate old\n// values (e.g. \"683x2=1366\") map to the correct distinct new values.\nconst table = context.document.body.tables.getFirst();\n\nconst answers = [\n    \"757\u00d75=3785\", \"223\u00d74=892\",  \"329\u00d79=2961\", \"932\u00d75=4660\", \"360\u00d78=2880\",\n    \"937\u00d74=3748\", \"762\u00d73=2286\", \"923\u00d78=7384\", \"582\u00d76=3492\", \"959\u00d73=2877\",\n    \"383\u00d72=766\",  \"518\u00d78=4144\", \"106\u00d79=954\",  \"158\u00d74=632\",  \"258\u00d72=516\",\n    \"569\u00d72=1138\", \"673\u00d77=4711\", \"330\u00d78=2640\", \"172\u00d77=1204\", \"818\u00d74=3272\",\n    \"986\u00d76=5916\", \"453\u00d76=2718\", \"969\u00d75=4845\", \"440\u00d78=3520\", \"149\u00d76=894\"\n];\n\nconst dataRows = [0, 4, 9, 14, 19];\nlet i = 0;\nfor (const row of dataRows) {\n    for (let col = 0; col < 5; col++) {\n        const cell = table.getCell(row, col);\n        const cellParagraph = cell.body.paragraphs.getFirst();\n        const cellRange = cellParagraph.getRange();\n        cellRange.insertText(answers[i], Word.InsertLocation.replace);\n        i++;\n    }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line at the top of the document.\n$dateRange = $d.Paragraphs.Item(1).Range\n$dateRange.MoveEnd(1, -1) | Out-Null\n$dateRange.Text = \"2026-02-14 Saturday\"\n\n# Update each three-digit x one-digit multiplication answer in the table.\n# The table has data only in rows 1, 5, 10, 15, 20 (5 columns each);\n# replacements are applied in row-major order so duplicate old values\n# (e.g. \"683x2=1366\") map to the correct distinct new values.\n$t = $d.Tables.Item(1)\n\n$answers = @(\n    \"757\u00d75=3785\", \"223\u00d74=892\",  \"329\u00d79=2961\", \"932\u00d75=4660\", \"360\u00d78=2880\",\n    \"937\u00d74=3748\", \"762\u00d73=2286\", \"923\u00d78=7384\", \"582\u00d76=3492\", \"959\u00d73=2877\",\n    \"383\u00d72=766\",  \"518\u00d78=4144\", \"106\u00d79=954\",  \"158\u00d74=632\",  \"258\u00d72=516\",\n    \"569\u00d72=1138\", \"673\u00d77=4711\", \"330\u00d78=2640\", \"172\u00d77=1204\", \"818\u00d74=3272\",\n    \"986\u00d76=5916\", \"453\u00d76=2718\", \"969\u00d75=4845\", \"440\u00d78=3520\", \"149\u00d76=894\"\n)\n\n$dataRows = @(1, 5, 10, 15, 20)\n$i = 0\nforeach ($row in $dataRows) {\n    for ($col = 1; $col -le 5; $col++) {\n        $cellRange = $t.Cell($row, $col).Range\n        $cellRange.MoveEnd(1, -1) | Out-Null\n        $cellRange.Text = $answers[$i]\n        $i++\n    }\n}\n"}
